$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 5.29586935043335
$ws.Range("B1").Value = 4.155534744262695
$ws.Range("C1").Value = 2.97716212272644
$ws.Range("D1").Value = 2.32318377494812
$ws.Range("E1").Value = 1.564671874046326
